# New weekly price observation for "Rabanito" at Vega Central Mapocho de
# Santiago: insert a row at position 64 (shifting the existing rows 64-163
# down to 65-164) and fill it in with the new week's figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value = 9
$ws.Cells.Item(64, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(64, 3).Value = "Metropolitana"
$ws.Cells.Item(64, 4).Value = 44477
$ws.Cells.Item(64, 5).Value = 13
$ws.Cells.Item(64, 6).Value = 300000001
$ws.Cells.Item(64, 7).Value = "Rabanito"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 7900
$ws.Cells.Item(64, 11).Value = 3500
$ws.Cells.Item(64, 12).Value = 4000
$ws.Cells.Item(64, 13).Value = 3747
$ws.Cells.Item(64, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(64, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(64, 16).Value = 37
$ws.Cells.Item(64, 17).Value = 100
$ws.Cells.Item(64, 18).Value = "Hortaliza"
